$d = $word.ActiveDocument

# Locate the paragraph that holds the Word field (fldChar begin/instrText/fldChar end)
# that encodes the m2doc query "m:'doc.html'.fromHTMLURI()". We find it generically
# (rather than hard-coding a paragraph index) by looking for a paragraph whose Range
# reports at least one Field.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
    }
}

if ($target -ne $null) {
    $xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>{</w:t></w:r>
<w:r><w:t>m</w:t></w:r>
<w:r><w:t>:</w:t></w:r>
<w:r><w:t>'</w:t></w:r>
<w:r><w:t>doc.html</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t>'.fromHTMLURI()</w:t></w:r>
<w:r><w:t xml:space="preserve">}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    [void]$target.Range.InsertXML($xml)
}
